$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")
$ws.Range("C2").Value = "ABF Freight, Ceva, FC Test Carrier, UPS, FragilePAK"
